$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Work breakdown structure update: split the "Upload sprints to github" (row 8)
# and "Take metrics" (row 9) effort across days H (4), I (5), J (6), K (7).

# Row 8 - "Upload sprints to github (Todo agente)"
$ws.Range("H8").Value = 0.5
$ws.Range("I8").Value = 0.5
$ws.Range("J8").Value = 0.5
$ws.Range("K8").Value = 0.5

# Row 9 - "Take metrics (Toda agente)"
$ws.Range("H9").ClearContents()
$ws.Range("I9").Value = 0.25
$ws.Range("J9").Value = 0.5
$ws.Range("K9").Value = 0.25

# Move the active selection to L9, matching the saved view state.
$ws.Range("L9").Select()
